$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 120.63636
$ws.Range("I2").Value = 83.57143000000001
$ws.Range("J2").Value = 185.5
$ws.Range("K2").Value = 83.57143000000001
$ws.Range("L2").Value = 185.5
$ws.Range("M2").Value = 29.42856999999999
$ws.Range("N2").Value = -411.5

$ws.Range("H18").Value = 595
$ws.Range("I18").Value = 595
$ws.Range("K18").Value = 595
$ws.Range("M18").Value = -311

$ws.Range("H43").Value = 2227.5
$ws.Range("I43").Value = 3075
$ws.Range("K43").Value = 3075
$ws.Range("M43").Value = -3006

$ws.Range("H64").Value = 3489.2666
$ws.Range("I64").Value = 3226.9
$ws.Range("J64").Value = 4014
$ws.Range("K64").Value = 3226.9
$ws.Range("L64").Value = 4014
$ws.Range("M64").Value = -2978.9
$ws.Range("N64").Value = -4510

$ws.Range("H67").Value = 3489.2666
$ws.Range("I67").Value = 3226.9
$ws.Range("J67").Value = 4014
$ws.Range("K67").Value = 3226.9
$ws.Range("L67").Value = 4014
$ws.Range("M67").Value = -2368.9
$ws.Range("N67").Value = -5730

$ws.Range("H76").Value = 7266.8823
$ws.Range("I76").Value = 8898.947
$ws.Range("J76").Value = 5199.6
$ws.Range("K76").Value = 8898.947
$ws.Range("L76").Value = 5199.6
$ws.Range("M76").Value = -8583.947
$ws.Range("N76").Value = -5829.6

$ws.Range("H79").Value = 7266.8823
$ws.Range("I79").Value = 8898.947
$ws.Range("J79").Value = 5199.6
$ws.Range("K79").Value = 8898.947
$ws.Range("L79").Value = 5199.6
$ws.Range("M79").Value = -7806.947
$ws.Range("N79").Value = -7383.6

$ws.Range("H137").Value = 18030.953
$ws.Range("I137").Value = 21737.53
$ws.Range("J137").Value = 2278
$ws.Range("K137").Value = 65212.59
$ws.Range("L137").Value = 6834
$ws.Range("M137").Value = -62662.59
$ws.Range("N137").Value = -11934

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7561.047
$ws.Range("I32").Value = 4831.1406
$ws.Range("J32").Value = 21405.572
$ws.Range("K32").Value = 4831.1406
$ws.Range("L32").Value = 21405.572
$ws.Range("M32").Value = -4544.1406
$ws.Range("N32").Value = -21979.572

$ws.Range("H45").Value = 1488.1
$ws.Range("I45").Value = 1251.091
$ws.Range("J45").Value = 1777.7778
$ws.Range("K45").Value = 1251.091
$ws.Range("L45").Value = 1777.7778
$ws.Range("M45").Value = -874.0909999999999
$ws.Range("N45").Value = -2531.7778

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1757.6957
$ws.Range("I20").Value = 1700.5385
$ws.Range("J20").Value = 1832
$ws.Range("K20").Value = 1700.5385
$ws.Range("L20").Value = 1832
$ws.Range("M20").Value = -1453.5385
$ws.Range("N20").Value = -2326

$ws.Range("H76").Value = 26000
$ws.Range("J76").Value = 26000
$ws.Range("L76").Value = 26000
$ws.Range("N76").Value = -26630

$ws.Range("H79").Value = 26000
$ws.Range("J79").Value = 26000
$ws.Range("L79").Value = 26000
$ws.Range("N79").Value = -28184

$ws.Range("H107").Value = 1881.7273
$ws.Range("I107").Value = 1674.4814
$ws.Range("J107").Value = 2814.3333
$ws.Range("K107").Value = 1674.4814
$ws.Range("L107").Value = 2814.3333
$ws.Range("M107").Value = 245.5186000000001
$ws.Range("N107").Value = -6654.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1300
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -713
$ws.Range("N16").Value = -2074

$ws.Range("H31").Value = 2348.541
$ws.Range("I31").Value = 1383.325
$ws.Range("J31").Value = 4187.048
$ws.Range("K31").Value = 1383.325
$ws.Range("L31").Value = 4187.048
$ws.Range("M31").Value = -1088.325
$ws.Range("N31").Value = -4777.048

$ws.Range("H34").Value = 2348.541
$ws.Range("I34").Value = 1383.325
$ws.Range("J34").Value = 4187.048
$ws.Range("K34").Value = 1383.325
$ws.Range("L34").Value = 4187.048
$ws.Range("M34").Value = -1181.325
$ws.Range("N34").Value = -4591.048

$ws.Range("H58").Value = 1417.7446
$ws.Range("I58").Value = 895.80646
$ws.Range("J58").Value = 2429
$ws.Range("K58").Value = 895.80646
$ws.Range("L58").Value = 2429
$ws.Range("M58").Value = -692.80646
$ws.Range("N58").Value = -2835

$ws.Range("H60").Value = 9997
$ws.Range("J60").Value = 20000
$ws.Range("L60").Value = 20000
$ws.Range("N60").Value = -21022

$ws.Range("H105").Value = 698.5714
$ws.Range("I105").Value = 613.8461
$ws.Range("K105").Value = 613.8461
$ws.Range("M105").Value = 1133.1539

$ws.Range("H113").Value = 1300
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 1170
$ws.Range("N113").Value = -5840

$ws.Range("H134").Value = 2287.1785
$ws.Range("I134").Value = 1403.2
$ws.Range("J134").Value = 2778.2778
$ws.Range("K134").Value = 4209.6
$ws.Range("L134").Value = 8334.8334
$ws.Range("M134").Value = -1674.6
$ws.Range("N134").Value = -13404.8334

$ws.Range("H136").Value = 1417.7446
$ws.Range("I136").Value = 895.80646
$ws.Range("J136").Value = 2429
$ws.Range("K136").Value = 2687.41938
$ws.Range("L136").Value = 7287
$ws.Range("M136").Value = -137.4193800000003
$ws.Range("N136").Value = -12387

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 4254.5454
$ws.Range("I56").Value = 4254.5454
$ws.Range("K56").Value = 4254.5454
$ws.Range("M56").Value = -3724.5454

$ws.Range("H60").Value = 155
$ws.Range("I60").Value = 97
$ws.Range("J60").Value = 300
$ws.Range("K60").Value = 291
$ws.Range("L60").Value = 900
$ws.Range("M60").Value = -40
$ws.Range("N60").Value = -1402

$ws.Range("H92").Value = 602.3889
$ws.Range("I92").Value = 373.1111
$ws.Range("J92").Value = 831.6667
$ws.Range("K92").Value = 1119.3333
$ws.Range("L92").Value = 2495.0001
$ws.Range("M92").Value = 128.6667
$ws.Range("N92").Value = -4991.0001

$ws.Range("H107").Value = 1192.3158
$ws.Range("J107").Value = 1670.7693
$ws.Range("L107").Value = 5012.3079
$ws.Range("N107").Value = -8852.3079

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5111.5386
$ws.Range("I80").Value = 5378.2607
$ws.Range("J80").Value = 3066.6667
$ws.Range("K80").Value = 5378.2607
$ws.Range("L80").Value = 3066.6667
$ws.Range("M80").Value = -4380.2607
$ws.Range("N80").Value = -5062.6667

$ws.Range("H83").Value = 5111.5386
$ws.Range("I83").Value = 5378.2607
$ws.Range("J83").Value = 3066.6667
$ws.Range("K83").Value = 26891.3035
$ws.Range("L83").Value = 15333.3335
$ws.Range("M83").Value = -21899.3035
$ws.Range("N83").Value = -25317.3335

$ws.Range("H122").Value = 949.8333
$ws.Range("I122").Value = 737.5
$ws.Range("J122").Value = 1374.5
$ws.Range("K122").Value = 2212.5
$ws.Range("L122").Value = 4123.5
$ws.Range("M122").Value = 237.5
$ws.Range("N122").Value = -9023.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1492.1389
$ws.Range("I16").Value = 1462.8667
$ws.Range("J16").Value = 1638.5
$ws.Range("K16").Value = 1462.8667
$ws.Range("L16").Value = 1638.5
$ws.Range("M16").Value = -1292.8667
$ws.Range("N16").Value = -1978.5

$ws.Range("H61").Value = 2815.8462
$ws.Range("I61").Value = 3040
$ws.Range("J61").Value = 2068.6667
$ws.Range("K61").Value = 3040
$ws.Range("L61").Value = 2068.6667
$ws.Range("M61").Value = -2838
$ws.Range("N61").Value = -2472.6667

$ws.Range("H68").Value = 2284.2856
$ws.Range("I68").Value = 1665
$ws.Range("J68").Value = 2748.75
$ws.Range("K68").Value = 1665
$ws.Range("L68").Value = 2748.75
$ws.Range("M68").Value = -916
$ws.Range("N68").Value = -4246.75

$ws.Range("H71").Value = 2284.2856
$ws.Range("I71").Value = 1665
$ws.Range("J71").Value = 2748.75
$ws.Range("K71").Value = 8325
$ws.Range("L71").Value = 13743.75
$ws.Range("M71").Value = -4581
$ws.Range("N71").Value = -21231.75

$ws.Range("H75").Value = 46174.4
$ws.Range("J75").Value = 46174.4
$ws.Range("L75").Value = 46174.4
$ws.Range("N75").Value = -48046.4

$ws.Range("H78").Value = 46174.4
$ws.Range("J78").Value = 46174.4
$ws.Range("L78").Value = 138523.2
$ws.Range("N78").Value = -147883.2

$ws.Range("H82").Value = 1847.4814
$ws.Range("I82").Value = 1346.2667
$ws.Range("J82").Value = 2474
$ws.Range("K82").Value = 1346.2667
$ws.Range("L82").Value = 2474
$ws.Range("M82").Value = -985.2666999999999
$ws.Range("N82").Value = -3196

$ws.Range("H85").Value = 1847.4814
$ws.Range("I85").Value = 1346.2667
$ws.Range("J85").Value = 2474
$ws.Range("K85").Value = 1346.2667
$ws.Range("L85").Value = 2474
$ws.Range("M85").Value = -98.2666999999999
$ws.Range("N85").Value = -4970

$ws.Range("H113").Value = 2815.8462
$ws.Range("I113").Value = 3040
$ws.Range("J113").Value = 2068.6667
$ws.Range("K113").Value = 3040
$ws.Range("L113").Value = 2068.6667
$ws.Range("M113").Value = -870
$ws.Range("N113").Value = -6408.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4123.273
$ws.Range("I62").Value = 3535.3333
$ws.Range("J62").Value = 4343.75
$ws.Range("K62").Value = 3535.3333
$ws.Range("L62").Value = 4343.75
$ws.Range("M62").Value = -2911.3333
$ws.Range("N62").Value = -5591.75

$ws.Range("H65").Value = 4123.273
$ws.Range("I65").Value = 3535.3333
$ws.Range("J65").Value = 4343.75
$ws.Range("K65").Value = 17676.6665
$ws.Range("L65").Value = 21718.75
$ws.Range("M65").Value = -14556.6665
$ws.Range("N65").Value = -27958.75

$ws.Range("H70").Value = 25000
$ws.Range("J70").Value = 25000
$ws.Range("L70").Value = 25000
$ws.Range("N70").Value = -25630

$ws.Range("H73").Value = 25000
$ws.Range("J73").Value = 25000
$ws.Range("L73").Value = 25000
$ws.Range("N73").Value = -27184

$ws.Range("H113").Value = 516.6667
$ws.Range("I113").Value = 533.3333
$ws.Range("J113").Value = 500
$ws.Range("K113").Value = 1599.9999
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 570.0001
$ws.Range("N113").Value = -5840

$ws.Range("H126").Value = 904.59576
$ws.Range("I126").Value = 690.1579
$ws.Range("J126").Value = 1810
$ws.Range("K126").Value = 2070.4737
$ws.Range("L126").Value = 5430
$ws.Range("M126").Value = 399.5263
$ws.Range("N126").Value = -10370
